$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.267.85'
$ws.Range('E2').Value = '  +2.28%  '
$ws.Range('D3').Value = '2.864.88'
$ws.Range('E3').Value = '  +7.77%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'194.10"
$ws.Range('E5').Value = '  +4.70%  '
$ws.Range('D6').Value = "'598.94"
$ws.Range('E6').Value = '  +2.73%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = "'0.552"
$ws.Range('E8').Value = '  +3.85%  '
$ws.Range('E9').Value = '  +1.77%  '
$ws.Range('D10').Value = '2.863.45'
$ws.Range('E10').Value = '  +7.82%  '
$ws.Range('E11').Value = '  +10.92%  '
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').Value = "'4.90"
$ws.Range('E13').Value = '  +3.45%  '
$ws.Range('D14').Value = '3.385.01'
$ws.Range('D15').Value = '76.111.78'
$ws.Range('E15').Value = '  +2.50%  '
$ws.Range('E16').Value = '  +4.58%  '
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('D18').Value = '2.865.39'
$ws.Range('E18').Value = '  +7.79%  '
$ws.Range('D19').Value = "'9.11"
$ws.Range('E19').Value = '  -1.62%  '
$ws.Range('D20').Value = "'12.48"
$ws.Range('E20').Value = '  +5.22%  '
$ws.Range('D21').Value = "'383.04"
$ws.Range('E21').Value = '  +3.46%  '
$ws.Range('D22').Value = "'2.34"
$ws.Range('E22').Value = '  +4.44%  '
$ws.Range('D23').Value = "'4.15"
$ws.Range('E23').Value = '  +2.28%  '
$ws.Range('D24').Value = "'71.98"
$ws.Range('E24').Value = '  +3.93%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '3.006.58'
$ws.Range('E26').Value = '  +7.70%  '
$ws.Range('B27').Value = 'NEARProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D27').Value = "'4.22"
$ws.Range('E27').Value = '  +3.43%  '
$ws.Range('D28').Value = "'9.79"
$ws.Range('E28').Value = '  +5.37%  '
$ws.Range('E29').Value = '  +12.95%  '
$ws.Range('D30').Value = "'1.01"
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('D32').Value = "'518.16"
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('D33').Value = "'7.72"
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('E34').Value = '  +4.96%  '
$ws.Range('D35').Value = "'0.999"
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = "'166.90"
$ws.Range('E36').Value = '  +2.17%  '
$ws.Range('D37').Value = "'20.06"
$ws.Range('E37').Value = '  +5.04%  '
$ws.Range('E38').Value = '  +1.80%  '
$ws.Range('D39').Value = "'19.46"
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').Value = "'186.34"
$ws.Range('E40').Value = '  +10.18%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').Value = "'0.346"
$ws.Range('E42').Value = '  +6.13%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').Value = "'5.10"
$ws.Range('E43').Value = '  +3.30%  '
$ws.Range('D44').Value = "'1.69"
$ws.Range('E44').Value = '  +1.87%  '
$ws.Range('D45').Value = "'1.23"
$ws.Range('E45').Value = '  +4.80%  '
$ws.Range('D46').Value = "'40.40"
$ws.Range('E46').Value = '  +3.55%  '
$ws.Range('D47').Value = "'0.0892"
$ws.Range('E47').Value = '  +5.38%  '
$ws.Range('D48').Value = "'2.37"
$ws.Range('E48').Value = '  +2.09%  '
$ws.Range('D49').Value = "'0.579"
$ws.Range('E49').Value = '  +10.74%  '
$ws.Range('D50').Value = "'3.77"
$ws.Range('E50').Value = '  +4.15%  '
$ws.Range('D51').Value = "'0.667"
$ws.Range('E51').Value = '  +13.20%  '
